# Update the mahasiswa (student) template sheet: drop EMAIL / Alamat / Tanggal
# Lahir / No Hp / Jenis Kelamin columns, keep Nama Lengkap / NPM / PRODI ID,
# and refresh the roster with the new data pulled from hosting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old contents (columns A:H, rows 1:4) and hyperlinks -----
$ws.Hyperlinks.Delete()
$ws.Range("A1:H4").Clear()

# --- Header row -------------------------------------------------------------
$ws.Range("A1").Value = "Nama Lengkap"
$ws.Range("B1").Value = "NPM"
$ws.Range("C1").Value = "PRODI ID"

# --- Data rows ---------------------------------------------------------------
$data = @(
    @("Agung Hardiyanto",      13753001, "800d3121-f2f9-11eb-8c3d-0cc47abcfaa6"),
    @("Jution Candra Kirana",  13753030, "800d3121-f2f9-11eb-8c3d-0cc47abcfaa6"),
    @("Adriyan Lumban Tobing", 13753004, "800d3121-f2f9-11eb-8c3d-0cc47abcfaa6"),
    @("Ria Setya Handayani",   13755001, "8f5585b1-f2f9-11eb-8c3d-0cc47abcfaa6"),
    @("Ayu Futri Batara",      13755002, "8f5585b1-f2f9-11eb-8c3d-0cc47abcfaa6")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# --- Column widths / layout ---------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.42578125
$ws.Columns.Item(2).ColumnWidth = 9.140625
$ws.Columns.Item(3).ColumnWidth = 37

$ws.Range("D1:G1").EntireColumn.Delete()

# --- Selection matching the saved view ---------------------------------------
$ws.Range("G1:G1048576").Select()
